$wb = $excel.ActiveWorkbook

# Add a new worksheet named "Sheet1" at the end (after "edges sheet")
$newWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$newWs.Name = "Sheet1"

$values = @(
    "10_11:10.307819,123.886707:10.305099,123.886696:10.303706,123.886846:10.300539,123.886718:10.297221,123.888381",
    "9_10:10.315672,123.885634:10.312952,123.886160:10.312404,123.886138:10.310250,123.886482:10.307819,123.886707",
    "4_2:10.289747,123.870742:10.290338,123.876278",
    "1_2:10.271590,123.847311:10.279418,123.855304:10.282817,123.858286:10.286976,123.862127:10.288834,123.866633:10.289747,123.870742",
    "30_4:10.270092,123.856403:10.270419,123.856779,10.274620,123.859912:10.274737,123.860115:10.283119,123.866714:10.284808,123.868516:10.290338,123.876278",
    "2_3:10.289747,123.870742:10.292429,123.870533:10.293685,123.870329:10.295469,123.869868:10.297243,123.869889:10.298372,123.869825:10.299639,123.869932:10.300959,123.869686:10.301443,123.870442",
    "4_7:10.290338,123.876278:10.296134,123.884395:10.296858,123.886136:10.297136,123.887994",
    "4_5:10.290338,123.876278:10.290634,123.877695",
    "3_6:10.301443,123.870442:10.299543,123.878982",
    "6_7:10.299543,123.878982:10.297136,123.887994",
    "7_11:10.297136,123.887994:10.297221,123.888381",
    "6_9:10.299543,123.878982:10.302737,123.877949:10.305207,123.876532:10.306157,123.876168:10.307339,123.875974:10.309873,123.876082:10.312469,123.877991:10.313630,123.880030:10.314665,123.882540:10.315672,123.885634",
    "8_9:10.315672,123.885634:10.323510,123.883721:10.324186,123.883356:10.324861,123.883227:10.329267,123.881085",
    "11_12:10.297221,123.888381:10.297660,123.890921:10.298315,123.892499:10.298705,123.895224:10.299459,123.896921",
    "12_13:10.299459,123.896921:10.297723,123.897530:10.296225,123.898239",
    "5_13:10.290634,123.877695:10.290883,123.878455:10.292920,123.885900:10.293860,123.890964:10.294251,123.893518:10.296165,123.898294",
    "9_20:10.315672,123.885634:10.316238,123.886909:10.317199,123.889838:10.316027,123.890235:10.316094,123.890955",
    "10_17:10.307819,123.886707:10.307942,123.887156:10.308121,123.889398:10.309339,123.892629",
    "20_18:10.316094,123.890955:10.310479,123.892887",
    "28_17:10.310479,123.892887:10.310179,123.892563:10.309715,123.892456:10.309339,123.892629",
    "17_16:10.309339,123.892629:10.309086,123.893445",
    "18_19:10.310479,123.892887:10.310475,123.893432:10.310232,123.893722:10.309803,123.893917",
    "19_16:10.309803,123.893917:10.309086,123.893445",
    "12_16:10.299459,123.896921:10.309086,123.893445",
    "13_15:10.296165,123.898294:10.293563,123.904592:10.294419,123.905094",
    "15_14:10.294419,123.905094:10.293970,123.906062",
    "12_24:10.299459,123.896921:10.301418,123.900320:10.302073,123.902133:10.302358,123.902444:10.308790,123.907135",
    "15_24:10.294419,123.905094:10.295623,123.905502:10.29727,123.905888:10.301313,123.906586:10.302315,123.906532:10.305018,123.906167:10.306358,123.906124:10.307393,123.906457:10.307994,123.906779:10.308438,123.90694:10.308790,123.907135",
    "20_21:10.316094,123.890955:10.317494,123.894623:10.318459,123.900011",
    "21_34:10.318459,123.900011:10.319933,123.899344:10.322413,123.898764:10.325326,123.898142:10.326941,123.897927:10.328641,123.897423:10.331131,123.898099:10.330998,123.898165",
    "21_22:10.318459,123.900011:10.31761,123.900459:10.315461,123.901684",
    "21_31:10.318459,123.900011:10.320429,123.903689:10.320317,123.903873",
    "22_31:10.315461,123.901684:10.317737,123.902627:10.318645,123.903056:10.320317,123.903873",
    "22_23:10.315461,123.901684:10.313156,123.902755:10.310732,123.903873",
    "19_23:10.309803,123.893917:10.311193,123.898206:10.311478,123.89958:10.311267,123.902444:10.310732,123.903873",
    "23_24:10.310732,123.903873:10.308790,123.907135",
    "24_25:10.308790,123.907135:10.307101,123.909924",
    "25_26:10.307101,123.909924:10.303597,123.912756",
    "26_27:10.303597,123.912756:10.308832,123.919151",
    "27_28:10.308832,123.919151:10.312294,123.916233",
    "28_25:10.312294,123.916233:10.307101,123.909924",
    "28_29:10.312294,123.916233:10.313768,123.914986:10.314194,123.914688",
    "29_24:10.314194,123.914688:10.311953,123.909944:10.310538,123.90857:10.308790,123.907135",
    "31_30:10.320317,123.903873:10.325932,123.907650",
    "34_30:10.330998,123.898165:10.3271,123.906081:10.325932,123.907650",
    "29_30:10.314194,123.914688:10.325932,123.907650",
    "14_26:10.293970,123.906062:10.293976,123.906028:10.301513,123.910201:10.303597,123.912756",
    "14_5:10.293970,123.906062:10.291569,123.90371:10.291231,123.903442:10.290197,123.899376:10.288677,123.897766:10.287917,123.889645:10.2819,123.883529:10.284919,123.880407:10.286798,123.879828:10.289342,123.87854:10.290634,123.877695",
    "30_32:10.325932,123.907650:10.326202,123.907905:10.331775,123.910276:10.333147,123.910523:10.339005,123.911585:10.339945,123.911832:10.349032,123.913334:10.350489,123.913549:10.351481,123.913817:10.352452,123.913988:10.354911,123.914739:10.355734,123.915094:10.357792,123.915169:10.358626,123.915705:10.360948,123.91534:10.362204,123.91534:10.364388,123.914654:10.36538,123.913838:10.367523,123.913913:10.371569:123.923957"
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 1
    $newWs.Cells.Item($row, 1).Value = $values[$i]
}

# Set view properties for the new sheet
$newWs.Application.ActiveWindow.ScrollRow = 43
$newWs.Range("A1:A49").Select()

# Re-select the original sheet and set its view
$ws1 = $wb.Worksheets.Item("edges sheet")
$ws1.Activate()
$ws1.Application.ActiveWindow.ScrollRow = 38
$ws1.Range("D51").Select()
